$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column date values for rows 2-6
# from 45221 (2023-10-22) to 45224 (2023-10-25), keeping existing formatting.
$ws.Range("C2:C6").Value = 45224
